$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: move the active selection from E9 to B7 ---
$ws.Range("B7").Select() | Out-Null

# --- Column widths: columns D and E were re-sized (best-fit to their
#     content) ---
# The host's ColumnWidth setter pads the value by a fixed 5/6 pt before it
# lands in the <col width="..."> attribute, so back that padding out to
# land exactly on the target stored widths (57.5 and 29.5).
$pad = 5/6
$ws.Columns.Item(4).ColumnWidth = 57.5 - $pad
$ws.Columns.Item(5).ColumnWidth = 29.5 - $pad

# --- Row heights: rows 9-22 (the blank template rows below the table)
#     grew from 18.3pt to 21.3pt, matching the single-line height of the
#     font already applied to those rows ---
for ($r = 9; $r -le 22; $r++) {
    $ws.Rows.Item($r).RowHeight = 21.3
}
